# Update weekly Fruta/Hortaliza price data.
# The dataset rotates: row 2 <- row 3, row 3 <- row 4, row 4 <- (old) row 2
# for the columns: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg).
#
# Note: use Value2 (not Value) to read/write, since Value returns a Variant
# wrapper object in this runtime that doesn't coerce cleanly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values before overwriting anything.
$origD2 = $ws.Range("D2").Value2
$origJ2 = $ws.Range("J2").Value2
$origK2 = $ws.Range("K2").Value2
$origL2 = $ws.Range("L2").Value2
$origM2 = $ws.Range("M2").Value2
$origP2 = $ws.Range("P2").Value2

$origD3 = $ws.Range("D3").Value2
$origJ3 = $ws.Range("J3").Value2
$origK3 = $ws.Range("K3").Value2
$origL3 = $ws.Range("L3").Value2
$origM3 = $ws.Range("M3").Value2
$origP3 = $ws.Range("P3").Value2

$origD4 = $ws.Range("D4").Value2
$origJ4 = $ws.Range("J4").Value2
$origK4 = $ws.Range("K4").Value2
$origL4 = $ws.Range("L4").Value2
$origM4 = $ws.Range("M4").Value2
$origP4 = $ws.Range("P4").Value2

# Row 2 gets old row 3's values
$ws.Range("D2").Value2 = $origD3
$ws.Range("J2").Value2 = $origJ3
$ws.Range("K2").Value2 = $origK3
$ws.Range("L2").Value2 = $origL3
$ws.Range("M2").Value2 = $origM3
$ws.Range("P2").Value2 = $origP3

# Row 3 gets old row 4's values
$ws.Range("D3").Value2 = $origD4
$ws.Range("J3").Value2 = $origJ4
$ws.Range("K3").Value2 = $origK4
$ws.Range("L3").Value2 = $origL4
$ws.Range("M3").Value2 = $origM4
$ws.Range("P3").Value2 = $origP4

# Row 4 gets old row 2's values
$ws.Range("D4").Value2 = $origD2
$ws.Range("J4").Value2 = $origJ2
$ws.Range("K4").Value2 = $origK2
$ws.Range("L4").Value2 = $origL2
$ws.Range("M4").Value2 = $origM2
$ws.Range("P4").Value2 = $origP2
